# Applies the "Update Excel data - 2024-11-22 04:11:09" refresh
# across all three sheets of the crypto live-data workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Top 50 Cryptocurrencies ---
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

$ws1.Cells.Item(2, 1).Value2 = "Bitcoin"
$ws1.Cells.Item(2, 2).Value2 = "btc"
$ws1.Cells.Item(2, 3).Value2 = 98936
$ws1.Cells.Item(2, 4).Value2 = 1957410186256
$ws1.Cells.Item(2, 5).Value2 = 115913150141
$ws1.Cells.Item(2, 6).Value2 = 2.57984

$ws1.Cells.Item(3, 1).Value2 = "Ethereum"
$ws1.Cells.Item(3, 2).Value2 = "eth"
$ws1.Cells.Item(3, 3).Value2 = 3371.17
$ws1.Cells.Item(3, 4).Value2 = 406015687428
$ws1.Cells.Item(3, 5).Value2 = 57864085710
$ws1.Cells.Item(3, 6).Value2 = 8.711589999999999

$ws1.Cells.Item(4, 1).Value2 = "Tether"
$ws1.Cells.Item(4, 2).Value2 = "usdt"
$ws1.Cells.Item(4, 3).Value2 = 1.002
$ws1.Cells.Item(4, 4).Value2 = 130920758489
$ws1.Cells.Item(4, 5).Value2 = 180658001561
$ws1.Cells.Item(4, 6).Value2 = -0.33176

$ws1.Cells.Item(5, 1).Value2 = "Solana"
$ws1.Cells.Item(5, 2).Value2 = "sol"
$ws1.Cells.Item(5, 3).Value2 = 262.94
$ws1.Cells.Item(5, 4).Value2 = 124821488522
$ws1.Cells.Item(5, 5).Value2 = 15497913447
$ws1.Cells.Item(5, 6).Value2 = 10.74741

$ws1.Cells.Item(6, 1).Value2 = "BNB"
$ws1.Cells.Item(6, 2).Value2 = "bnb"
$ws1.Cells.Item(6, 3).Value2 = 636.66
$ws1.Cells.Item(6, 4).Value2 = 92924369009
$ws1.Cells.Item(6, 5).Value2 = 2560789290
$ws1.Cells.Item(6, 6).Value2 = 5.02412

$ws1.Cells.Item(7, 1).Value2 = "XRP"
$ws1.Cells.Item(7, 2).Value2 = "xrp"
$ws1.Cells.Item(7, 3).Value2 = 1.38
$ws1.Cells.Item(7, 4).Value2 = 78666166648
$ws1.Cells.Item(7, 5).Value2 = 17004059380
$ws1.Cells.Item(7, 6).Value2 = 25.01952

$ws1.Cells.Item(8, 1).Value2 = "Dogecoin"
$ws1.Cells.Item(8, 2).Value2 = "doge"
$ws1.Cells.Item(8, 3).Value2 = 0.393894
$ws1.Cells.Item(8, 4).Value2 = 57853509205
$ws1.Cells.Item(8, 5).Value2 = 10366188434
$ws1.Cells.Item(8, 6).Value2 = 3.63818

$ws1.Cells.Item(9, 1).Value2 = "USDC"
$ws1.Cells.Item(9, 2).Value2 = "usdc"
$ws1.Cells.Item(9, 3).Value2 = 1
$ws1.Cells.Item(9, 4).Value2 = 38253314257
$ws1.Cells.Item(9, 5).Value2 = 17378138495
$ws1.Cells.Item(9, 6).Value2 = -0.29745

$ws1.Cells.Item(10, 1).Value2 = "Lido Staked Ether"
$ws1.Cells.Item(10, 2).Value2 = "steth"
$ws1.Cells.Item(10, 3).Value2 = 3367.64
$ws1.Cells.Item(10, 4).Value2 = 32985342102
$ws1.Cells.Item(10, 5).Value2 = 149669851
$ws1.Cells.Item(10, 6).Value2 = 8.403560000000001

$ws1.Cells.Item(11, 1).Value2 = "Cardano"
$ws1.Cells.Item(11, 2).Value2 = "ada"
$ws1.Cells.Item(11, 3).Value2 = 0.869827
$ws1.Cells.Item(11, 4).Value2 = 31086812023
$ws1.Cells.Item(11, 5).Value2 = 3350520127
$ws1.Cells.Item(11, 6).Value2 = 11.21959

$ws1.Cells.Item(12, 1).Value2 = "TRON"
$ws1.Cells.Item(12, 2).Value2 = "trx"
$ws1.Cells.Item(12, 3).Value2 = 0.200809
$ws1.Cells.Item(12, 4).Value2 = 17337544762
$ws1.Cells.Item(12, 5).Value2 = 1063110796
$ws1.Cells.Item(12, 6).Value2 = 2.15375

$ws1.Cells.Item(13, 1).Value2 = "Avalanche"
$ws1.Cells.Item(13, 2).Value2 = "avax"
$ws1.Cells.Item(13, 3).Value2 = 36.29
$ws1.Cells.Item(13, 4).Value2 = 14842978811
$ws1.Cells.Item(13, 5).Value2 = 1056089658
$ws1.Cells.Item(13, 6).Value2 = 8.61652

$ws1.Cells.Item(14, 1).Value2 = "Shiba Inu"
$ws1.Cells.Item(14, 2).Value2 = "shib"
$ws1.Cells.Item(14, 3).Value2 = 0.00002499
$ws1.Cells.Item(14, 4).Value2 = 14710554787
$ws1.Cells.Item(14, 5).Value2 = 1616603700
$ws1.Cells.Item(14, 6).Value2 = 5.48248

$ws1.Cells.Item(15, 1).Value2 = "Wrapped Bitcoin"
$ws1.Cells.Item(15, 2).Value2 = "wbtc"
$ws1.Cells.Item(15, 3).Value2 = 98586
$ws1.Cells.Item(15, 4).Value2 = 14398896989
$ws1.Cells.Item(15, 5).Value2 = 921826657
$ws1.Cells.Item(15, 6).Value2 = 2.38183

$ws1.Cells.Item(16, 1).Value2 = "Wrapped stETH"
$ws1.Cells.Item(16, 2).Value2 = "wsteth"
$ws1.Cells.Item(16, 3).Value2 = 3988.96
$ws1.Cells.Item(16, 4).Value2 = 14386360829
$ws1.Cells.Item(16, 5).Value2 = 167163549
$ws1.Cells.Item(16, 6).Value2 = 8.905950000000001

$ws1.Cells.Item(17, 1).Value2 = "Toncoin"
$ws1.Cells.Item(17, 2).Value2 = "ton"
$ws1.Cells.Item(17, 3).Value2 = 5.55
$ws1.Cells.Item(17, 4).Value2 = 14144705757
$ws1.Cells.Item(17, 5).Value2 = 621432219
$ws1.Cells.Item(17, 6).Value2 = 4.44313

$ws1.Cells.Item(18, 1).Value2 = "Sui"
$ws1.Cells.Item(18, 2).Value2 = "sui"
$ws1.Cells.Item(18, 3).Value2 = 3.6
$ws1.Cells.Item(18, 4).Value2 = 10231629493
$ws1.Cells.Item(18, 5).Value2 = 2287578574
$ws1.Cells.Item(18, 6).Value2 = 2.74334

$ws1.Cells.Item(19, 1).Value2 = "Bitcoin Cash"
$ws1.Cells.Item(19, 2).Value2 = "bch"
$ws1.Cells.Item(19, 3).Value2 = 494.5
$ws1.Cells.Item(19, 4).Value2 = 9798830623
$ws1.Cells.Item(19, 5).Value2 = 2350395615
$ws1.Cells.Item(19, 6).Value2 = 7.71998

$ws1.Cells.Item(20, 1).Value2 = "WETH"
$ws1.Cells.Item(20, 2).Value2 = "weth"
$ws1.Cells.Item(20, 3).Value2 = 3373.33
$ws1.Cells.Item(20, 4).Value2 = 9708765975
$ws1.Cells.Item(20, 5).Value2 = 519332627
$ws1.Cells.Item(20, 6).Value2 = 8.80829

$ws1.Cells.Item(21, 1).Value2 = "Chainlink"
$ws1.Cells.Item(21, 2).Value2 = "link"
$ws1.Cells.Item(21, 3).Value2 = 15.12
$ws1.Cells.Item(21, 4).Value2 = 9500220377
$ws1.Cells.Item(21, 5).Value2 = 1223424068
$ws1.Cells.Item(21, 6).Value2 = 6.34317

$ws1.Cells.Item(22, 1).Value2 = "Pepe"
$ws1.Cells.Item(22, 2).Value2 = "pepe"
$ws1.Cells.Item(22, 3).Value2 = 0.00002142
$ws1.Cells.Item(22, 4).Value2 = 9014421360
$ws1.Cells.Item(22, 5).Value2 = 7032492661
$ws1.Cells.Item(22, 6).Value2 = 13.29684

$ws1.Cells.Item(23, 1).Value2 = "Polkadot"
$ws1.Cells.Item(23, 2).Value2 = "dot"
$ws1.Cells.Item(23, 3).Value2 = 6.17
$ws1.Cells.Item(23, 4).Value2 = 8881549580
$ws1.Cells.Item(23, 5).Value2 = 814547946
$ws1.Cells.Item(23, 6).Value2 = 9.85849

$ws1.Cells.Item(24, 1).Value2 = "Stellar"
$ws1.Cells.Item(24, 2).Value2 = "xlm"
$ws1.Cells.Item(24, 3).Value2 = 0.283117
$ws1.Cells.Item(24, 4).Value2 = 8501860012
$ws1.Cells.Item(24, 5).Value2 = 2311728073
$ws1.Cells.Item(24, 6).Value2 = 17.79302

$ws1.Cells.Item(25, 1).Value2 = "LEO Token"
$ws1.Cells.Item(25, 2).Value2 = "leo"
$ws1.Cells.Item(25, 3).Value2 = 8.76
$ws1.Cells.Item(25, 4).Value2 = 8078033416
$ws1.Cells.Item(25, 5).Value2 = 3517774
$ws1.Cells.Item(25, 6).Value2 = 2.67332

$ws1.Cells.Item(26, 1).Value2 = "NEAR Protocol"
$ws1.Cells.Item(26, 2).Value2 = "near"
$ws1.Cells.Item(26, 3).Value2 = 5.78
$ws1.Cells.Item(26, 4).Value2 = 7035895518
$ws1.Cells.Item(26, 5).Value2 = 968924810
$ws1.Cells.Item(26, 6).Value2 = 5.61817

$ws1.Cells.Item(27, 1).Value2 = "Litecoin"
$ws1.Cells.Item(27, 2).Value2 = "ltc"
$ws1.Cells.Item(27, 3).Value2 = 90.19
$ws1.Cells.Item(27, 4).Value2 = 6784076065
$ws1.Cells.Item(27, 5).Value2 = 1460843275
$ws1.Cells.Item(27, 6).Value2 = 6.54217

$ws1.Cells.Item(28, 1).Value2 = "Aptos"
$ws1.Cells.Item(28, 2).Value2 = "apt"
$ws1.Cells.Item(28, 3).Value2 = 12.1
$ws1.Cells.Item(28, 4).Value2 = 6453451196
$ws1.Cells.Item(28, 5).Value2 = 895334881
$ws1.Cells.Item(28, 6).Value2 = 5.04982

$ws1.Cells.Item(29, 1).Value2 = "Wrapped eETH"
$ws1.Cells.Item(29, 2).Value2 = "weeth"
$ws1.Cells.Item(29, 3).Value2 = 3544.81
$ws1.Cells.Item(29, 4).Value2 = 6101083612
$ws1.Cells.Item(29, 5).Value2 = 96814193
$ws1.Cells.Item(29, 6).Value2 = 9.20321

$ws1.Cells.Item(30, 1).Value2 = "Uniswap"
$ws1.Cells.Item(30, 2).Value2 = "uni"
$ws1.Cells.Item(30, 3).Value2 = 9.35
$ws1.Cells.Item(30, 4).Value2 = 5609709607
$ws1.Cells.Item(30, 5).Value2 = 847563217
$ws1.Cells.Item(30, 6).Value2 = 7.90214

$ws1.Cells.Item(31, 1).Value2 = "Cronos"
$ws1.Cells.Item(31, 2).Value2 = "cro"
$ws1.Cells.Item(31, 3).Value2 = 0.196843
$ws1.Cells.Item(31, 4).Value2 = 5336019451
$ws1.Cells.Item(31, 5).Value2 = 115416917
$ws1.Cells.Item(31, 6).Value2 = 10.11105

$ws1.Cells.Item(32, 1).Value2 = "USDS"
$ws1.Cells.Item(32, 2).Value2 = "usds"
$ws1.Cells.Item(32, 3).Value2 = 1.008
$ws1.Cells.Item(32, 4).Value2 = 5286602021
$ws1.Cells.Item(32, 5).Value2 = 16038564
$ws1.Cells.Item(32, 6).Value2 = 0.24725

$ws1.Cells.Item(33, 1).Value2 = "Hedera"
$ws1.Cells.Item(33, 2).Value2 = "hbar"
$ws1.Cells.Item(33, 3).Value2 = 0.129052
$ws1.Cells.Item(33, 4).Value2 = 4902210451
$ws1.Cells.Item(33, 5).Value2 = 842541829
$ws1.Cells.Item(33, 6).Value2 = 1.16109

$ws1.Cells.Item(34, 1).Value2 = "Internet Computer"
$ws1.Cells.Item(34, 2).Value2 = "icp"
$ws1.Cells.Item(34, 3).Value2 = 9.609999999999999
$ws1.Cells.Item(34, 4).Value2 = 4560151261
$ws1.Cells.Item(34, 5).Value2 = 272383770
$ws1.Cells.Item(34, 6).Value2 = 8.18369

$ws1.Cells.Item(35, 1).Value2 = "Ethereum Classic"
$ws1.Cells.Item(35, 2).Value2 = "etc"
$ws1.Cells.Item(35, 3).Value2 = 27.83
$ws1.Cells.Item(35, 4).Value2 = 4164131514
$ws1.Cells.Item(35, 5).Value2 = 919838112
$ws1.Cells.Item(35, 6).Value2 = 8.69515

$ws1.Cells.Item(36, 1).Value2 = "Bonk"
$ws1.Cells.Item(36, 2).Value2 = "bonk"
$ws1.Cells.Item(36, 3).Value2 = 0.00005245
$ws1.Cells.Item(36, 4).Value2 = 3920752548
$ws1.Cells.Item(36, 5).Value2 = 1818347499
$ws1.Cells.Item(36, 6).Value2 = 8.48047

$ws1.Cells.Item(37, 1).Value2 = "Render"
$ws1.Cells.Item(37, 2).Value2 = "render"
$ws1.Cells.Item(37, 3).Value2 = 7.42
$ws1.Cells.Item(37, 4).Value2 = 3842357119
$ws1.Cells.Item(37, 5).Value2 = 447271818
$ws1.Cells.Item(37, 6).Value2 = 2.09978

$ws1.Cells.Item(38, 1).Value2 = "Kaspa"
$ws1.Cells.Item(38, 2).Value2 = "kas"
$ws1.Cells.Item(38, 3).Value2 = 0.150282
$ws1.Cells.Item(38, 4).Value2 = 3786320788
$ws1.Cells.Item(38, 5).Value2 = 156046143
$ws1.Cells.Item(38, 6).Value2 = 1.04022

$ws1.Cells.Item(39, 1).Value2 = "Bittensor"
$ws1.Cells.Item(39, 2).Value2 = "tao"
$ws1.Cells.Item(39, 3).Value2 = 509.21
$ws1.Cells.Item(39, 4).Value2 = 3761174767
$ws1.Cells.Item(39, 5).Value2 = 268386450
$ws1.Cells.Item(39, 6).Value2 = 5.3129

$ws1.Cells.Item(40, 1).Value2 = "POL (ex-MATIC)"
$ws1.Cells.Item(40, 2).Value2 = "pol"
$ws1.Cells.Item(40, 3).Value2 = 0.469057
$ws1.Cells.Item(40, 4).Value2 = 3735128446
$ws1.Cells.Item(40, 5).Value2 = 448961333
$ws1.Cells.Item(40, 6).Value2 = 8.945919999999999

$ws1.Cells.Item(41, 1).Value2 = "Ethena USDe"
$ws1.Cells.Item(41, 2).Value2 = "usde"
$ws1.Cells.Item(41, 3).Value2 = 1.001
$ws1.Cells.Item(41, 4).Value2 = 3683475691
$ws1.Cells.Item(41, 5).Value2 = 241476759
$ws1.Cells.Item(41, 6).Value2 = -0.74746

$ws1.Cells.Item(42, 1).Value2 = "WhiteBIT Coin"
$ws1.Cells.Item(42, 2).Value2 = "wbt"
$ws1.Cells.Item(42, 3).Value2 = 24.8
$ws1.Cells.Item(42, 4).Value2 = 3573075283
$ws1.Cells.Item(42, 5).Value2 = 39356445
$ws1.Cells.Item(42, 6).Value2 = 2.5274

$ws1.Cells.Item(43, 1).Value2 = "MANTRA"
$ws1.Cells.Item(43, 2).Value2 = "om"
$ws1.Cells.Item(43, 3).Value2 = 3.9
$ws1.Cells.Item(43, 4).Value2 = 3506831676
$ws1.Cells.Item(43, 5).Value2 = 303461473
$ws1.Cells.Item(43, 6).Value2 = 6.54558

$ws1.Cells.Item(44, 1).Value2 = "Dai"
$ws1.Cells.Item(44, 2).Value2 = "dai"
$ws1.Cells.Item(44, 3).Value2 = 1.001
$ws1.Cells.Item(44, 4).Value2 = 3447798190
$ws1.Cells.Item(44, 5).Value2 = 190616435
$ws1.Cells.Item(44, 6).Value2 = -0.27302

$ws1.Cells.Item(45, 1).Value2 = "Artificial Superintelligence Alliance"
$ws1.Cells.Item(45, 2).Value2 = "fet"
$ws1.Cells.Item(45, 3).Value2 = 1.28
$ws1.Cells.Item(45, 4).Value2 = 3348828073
$ws1.Cells.Item(45, 5).Value2 = 499729195
$ws1.Cells.Item(45, 6).Value2 = 5.2388

$ws1.Cells.Item(46, 1).Value2 = "dogwifhat"
$ws1.Cells.Item(46, 2).Value2 = "wif"
$ws1.Cells.Item(46, 3).Value2 = 3.34
$ws1.Cells.Item(46, 4).Value2 = 3328766836
$ws1.Cells.Item(46, 5).Value2 = 1254384009
$ws1.Cells.Item(46, 6).Value2 = 7.44185

$ws1.Cells.Item(47, 1).Value2 = "Arbitrum"
$ws1.Cells.Item(47, 2).Value2 = "arb"
$ws1.Cells.Item(47, 3).Value2 = 0.772594
$ws1.Cells.Item(47, 4).Value2 = 3165596310
$ws1.Cells.Item(47, 5).Value2 = 1666852762
$ws1.Cells.Item(47, 6).Value2 = 14.65999

$ws1.Cells.Item(48, 1).Value2 = "Monero"
$ws1.Cells.Item(48, 2).Value2 = "xmr"
$ws1.Cells.Item(48, 3).Value2 = 160.43
$ws1.Cells.Item(48, 4).Value2 = 2960160649
$ws1.Cells.Item(48, 5).Value2 = 84063512
$ws1.Cells.Item(48, 6).Value2 = -0.91329

$ws1.Cells.Item(49, 1).Value2 = "Stacks"
$ws1.Cells.Item(49, 2).Value2 = "stx"
$ws1.Cells.Item(49, 3).Value2 = 1.95
$ws1.Cells.Item(49, 4).Value2 = 2922681927
$ws1.Cells.Item(49, 5).Value2 = 435064718
$ws1.Cells.Item(49, 6).Value2 = 4.65692

$ws1.Cells.Item(50, 1).Value2 = "OKB"
$ws1.Cells.Item(50, 2).Value2 = "okb"
$ws1.Cells.Item(50, 3).Value2 = 46.63
$ws1.Cells.Item(50, 4).Value2 = 2800201189
$ws1.Cells.Item(50, 5).Value2 = 20274712
$ws1.Cells.Item(50, 6).Value2 = 6.08207

$ws1.Cells.Item(51, 1).Value2 = "Filecoin"
$ws1.Cells.Item(51, 2).Value2 = "fil"
$ws1.Cells.Item(51, 3).Value2 = 4.66
$ws1.Cells.Item(51, 4).Value2 = 2795603507
$ws1.Cells.Item(51, 5).Value2 = 583510661
$ws1.Cells.Item(51, 6).Value2 = 9.148809999999999

# --- Sheet 2: Top 5 by Market Cap ---
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")

$ws2.Cells.Item(2, 2).Value2 = 1957410186256
$ws2.Cells.Item(3, 2).Value2 = 406015687428
$ws2.Cells.Item(4, 2).Value2 = 130920758489
$ws2.Cells.Item(5, 2).Value2 = 124821488522
$ws2.Cells.Item(6, 2).Value2 = 92924369009

# --- Sheet 3: Summary ---
$ws3 = $wb.Worksheets.Item("Summary")
# Leading "$" would otherwise auto-parse as currency; force literal text with
# the classic apostrophe prefix, then clear the resulting "quote prefix" style
# so the cell ends up with plain default formatting (no explicit style index).
$ws3.Cells.Item(2, 2).Value2 = "'" + '$4351.32'
$ws3.Cells.Item(2, 2).Style = "Normal"
$ws3.Cells.Item(3, 2).Value2 = "XRP (25.02%)"
$ws3.Cells.Item(4, 2).Value2 = "Monero (-0.91%)"
